$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the Status text for all data rows from OPTIMAL to TIME_LIMIT
$ws.Range("E2:E11").Value = "TIME_LIMIT"

# Corrected fixed recourse data: objective (B), gap (C), solve time (D)
$data = @(
    @(-1726.3126788491722, 8.882671417141983, 5918.626189953),
    @(-1740.0765152806462, 8.027488494657971, 5662.842745404),
    @(-1732.2668486303287, 9.519167611773694, 5567.290150216),
    @(-1731.213773381394,  10.14910946044649, 5563.251433826),
    @(-1743.143823798119,  7.643987046190622, 5581.206660045),
    @(-1719.7767196432865, 7.781098670675725, 5601.824704191),
    @(-1732.690555625326,  8.88929688644903,  5623.321135761),
    @(-1716.7979290041626, 9.388813050529913, 5614.556414993),
    @(-1730.9641896387025, 7.424981119889535, 5649.392891892),
    @(-1743.1617425024665, 8.753854342914902, 5625.161858239)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
